$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark ft_itoa.c (row 11), ft_split.c (row 25) and ft_substr.c (row 33) as
# "Working" (x) by copying the already-completed formatting (green fill)
# from a row that is already marked, then writing the "x" value into the
# Working column.
$ws.Range("C4:F4").Copy() | Out-Null
$ws.Range("C11:F11").PasteSpecial(-4122) | Out-Null
$ws.Range("C25:F25").PasteSpecial(-4122) | Out-Null
$ws.Range("C33:F33").PasteSpecial(-4122) | Out-Null

$ws.Range("F11").Value = "x"
$ws.Range("F25").Value = "x"
$ws.Range("F33").Value = "x"

# Restore the selection / active cell to match the saved workbook state.
$ws.Range("I12").Select()
